$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates
# NOTE: values that look like plain numbers (single decimal point) are
# prefixed with a leading apostrophe so Excel stores them as text
# (matching the original inlineStr/string cell type) instead of
# auto-converting them to a numeric value.
$ws.Range("D2").Value  = "57.138.30"
$ws.Range("D3").Value  = "3.036.30"
$ws.Range("D5").Value  = "'513.22"
$ws.Range("D6").Value  = "'140.03"
$ws.Range("D12").Value = "3.574.51"
$ws.Range("D14").Value = "'26.87"
$ws.Range("D16").Value = "57.080.02"
$ws.Range("D17").Value = "'6.17"
$ws.Range("D18").Value = "3.040.67"
$ws.Range("D19").Value = "'13.36"
$ws.Range("D21").Value = "'329.53"
$ws.Range("D23").Value = "'0.504"
$ws.Range("D24").Value = "'65.30"
$ws.Range("D25").Value = "3.175.41"
$ws.Range("D28").Value = "0.0₃0884"
$ws.Range("D29").Value = "'6.68"
$ws.Range("D31").Value = "'1.80"
$ws.Range("D32").Value = "'1.20"
$ws.Range("D33").Value = "'20.70"
$ws.Range("D35").Value = "'152.58"
$ws.Range("D38").Value = "'25.19"
$ws.Range("D43").Value = "'0.661"
$ws.Range("D45").Value = "2.193.51"
$ws.Range("D46").Value = "'6.08"
$ws.Range("D47").Value = "'0.946"
$ws.Range("D51").Value = "'0.0864"

# Column E (Volume 1h) updates
$ws.Range("E2").Value  = "  -0.46%  "
$ws.Range("E3").Value  = "  +0.79%  "
$ws.Range("E4").Value  = "  -0.01%  "
$ws.Range("E5").Value  = "  +1.03%  "
$ws.Range("E6").Value  = "  +0.46%  "
$ws.Range("E7").Value  = "  +0.02%  "
$ws.Range("E8").Value  = "  +0.79%  "
$ws.Range("E9").Value  = "  -5.21%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("E28").Value = "  -3.65%  "
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  -1.85%  "
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +1.78%  "
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  -3.41%  "

# Rows 48 and 49 swap (InjectiveProtocol <-> VeChain)
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0241"
$ws.Range("E48").Value = "  +0.76%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'20.07"
$ws.Range("E49").Value = "  +2.79%  "
